$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45180 -> 45181)
# for every data row from row 2 through row 99. Update the value in
# place, preserving existing formatting/style.
for ($row = 2; $row -le 99; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
